$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = 5
$ws.Range("F11").Value = 0
$ws.Range("F16").Value = -5
